$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Nppa", "Npr1", "ECs", 1, 0.3333333333333333, 0.110293, 0.330879, 0.1397687944778318, 0.1397687944778318, 3, 1, 55.59737033333334, 166.792111, 0.9022741440276701, 0.9022741440276703, 6.132000766174333, 55.188006895569, 0.126109769399265, 0.126109769399265)
$row3 = @("ECs", "Nppa", "Npr1", "FAPs", 1, 0.3333333333333333, 0.110293, 0.330879, 0.1397687944778318, 0.1397687944778318, 3, 1, 5.463788333333333, 16.391365, 0.08867028983655055, 0.08867028983655056, 0.6026176066483333, 5.423558459835, 0.01239333951645461, 0.01239333951645461)
$row4 = @("ECs", "Nppa", "Npr1", "MuSCs", 1, 0.3333333333333333, 0.110293, 0.330879, 0.1397687944778318, 0.1397687944778318, 3, 1, 0.4832879999999999, 1.449864, 0.007843145528366948, 0.00784314552836695, 0.05330328338399998, 0.4797295504559999, 0.001096226995414045, 0.001096226995414046)
$row5 = @("ECs", "Nppa", "Npr1", "Resolving-Mac", 1, 0.3333333333333333, 0.110293, 0.330879, 0.1397687944778318, 0.1397687944778318, 1, 0.3333333333333333, 0.07470833333333333, 0.224125, 0.001212420607412311, 0.001212420607412311, 0.008239806208333332, 0.074158255875, 0.0001694585666980992, 0.0001694585666980992)
$row6 = @("FAPs", "Nppa", "Npr1", "ECs", 1, 0.3333333333333333, 0.03729466666666666, 0.111884, 0.04726166302895539, 0.04726166302895539, 3, 1, 55.59737033333334, 166.792111, 0.9022741440276701, 0.9022741440276703, 2.073485394124889, 18.661368547124, 0.04264297655477491, 0.04264297655477491)
$row7 = @("FAPs", "Nppa", "Npr1", "FAPs", 1, 0.3333333333333333, 0.03729466666666666, 0.111884, 0.04726166302895539, 0.04726166302895539, 3, 1, 5.463788333333333, 16.391365, 0.08867028983655055, 0.08867028983655056, 0.2037701646288889, 1.83393148166, 0.00419070535893486, 0.00419070535893486)
$row8 = @("FAPs", "Nppa", "Npr1", "MuSCs", 1, 0.3333333333333333, 0.03729466666666666, 0.111884, 0.04726166302895539, 0.04726166302895539, 3, 1, 0.2948553333333333, 0.884566, 0.008161772060237412, 0.00816177206023741, 0.018024064864, 0.162216583776, 0.000370680101048737, 0.0003706801010487371)
$row9 = @("FAPs", "Nppa", "Npr1", "Resolving-Mac", 1, 0.3333333333333333, 0.03729466666666666, 0.111884, 0.04726166302895539, 0.04726166302895539, 1, 0.3333333333333333, 0.07470833333333333, 0.224125, 0.001212420607412311, 0.001212420607412311, 0.002786222388888889, 0.0250760015, 0.00005730101419688205, 0.00005730101419688205)
$row10 = @("MuSCs", "Nppa", "Npr1", "ECs", 2, 0.6666666666666666, 0.3940343333333334, 1.182103, 0.4993399740044802, 0.4993399740044802, 3, 1, 55.59737033333334, 166.792111, 0.9022741440276701, 0.9022741440276703, 21.90727275438145, 197.165454789433, 0.4505415476236914, 0.4505415476236915)
$row11 = @("MuSCs", "Nppa", "Npr1", "FAPs", 2, 0.6666666666666666, 0.3940343333333334, 1.182103, 0.4993399740044802, 0.4993399740044802, 3, 1, 5.463788333333333, 16.391365, 0.08867028983655055, 0.08867028983655056, 2.152920193399445, 19.376281740595, 0.04427662022195288, 0.04427662022195288)
$row12 = @("MuSCs", "Nppa", "Npr1", "MuSCs", 2, 0.6666666666666666, 0.3940343333333334, 1.182103, 0.4993399740044802, 0.4993399740044802, 3, 1, 0.4832879999999999, 1.449864, 0.007843145528366948, 0.00784314552836695, 0.190432064888, 1.713888583992, 0.003916396084248108, 0.003916396084248108)
$row13 = @("MuSCs", "Nppa", "Npr1", "Resolving-Mac", 2, 0.6666666666666666, 0.3940343333333334, 1.182103, 0.4993399740044802, 0.4993399740044802, 1, 0.3333333333333333, 0.07470833333333333, 0.224125, 0.001212420607412311, 0.001212420607412311, 0.02943764831944445, 0.264938834875, 0.0006054100745877594, 0.0006054100745877594)
$row14 = @("Resolving-Mac", "Nppa", "Npr1", "ECs", 1, 0.3333333333333333, 0.2474883333333333, 0.742465, 0.3136295684887327, 0.3136295684887327, 3, 1, 55.59737033333334, 166.792111, 0.9022741440276701, 0.9022741440276703, 13.75970052151278, 123.837304693615, 0.2829798504499388, 0.2829798504499389)
$row15 = @("Resolving-Mac", "Nppa", "Npr1", "FAPs", 1, 0.3333333333333333, 0.2474883333333333, 0.742465, 0.3136295684887327, 0.3136295684887327, 3, 1, 5.463788333333333, 16.391365, 0.08867028983655055, 0.08867028983655056, 1.352223868302778, 12.170014814725, 0.02780962473920821, 0.02780962473920821)
$row16 = @("Resolving-Mac", "Nppa", "Npr1", "MuSCs", 1, 0.3333333333333333, 0.2474883333333333, 0.742465, 0.3136295684887327, 0.3136295684887327, 3, 1, 0.4832879999999999, 1.449864, 0.007843145528366948, 0.00784314552836695, 0.11960814164, 1.07647327476, 0.00245984234765606, 0.00245984234765606)
$row17 = @("Resolving-Mac", "Nppa", "Npr1", "Resolving-Mac", 1, 0.3333333333333333, 0.2474883333333333, 0.742465, 0.3136295684887327, 0.3136295684887327, 1, 0.3333333333333333, 0.07470833333333333, 0.224125, 0.001212420607412311, 0.001212420607412311, 0.01848944090277778, 0.166404968125, 0.0003802509519295702, 0.0003802509519295702)

$allRows = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13, $row14, $row15, $row16, $row17)

$r = 2
foreach ($rowVals in $allRows) {
    $c = 1
    foreach ($v in $rowVals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
    $r = $r + 1
}
